{"js": "// Remove the resolved bug-tracker paragraphs from the body.\n// These are the exact paragraph texts that must be deleted entirely\n// (whole paragraph, including every run inside it).\nconst textsToRemove = [\n  \"Matricule doit etre string ou int : parce que si int 00012 devient 12 -> fait\",\n  \"Mettre automatique la barre / pour les dates ?? -> fait\",\n  \"Modifier les vue pour ne pas utilise dispose sur la vue\",\n  \"Est ce qu\\u2019on supprime tous les objet qui vont avec un genre on supprime les livraison d\\u2019un chantier si on supprime le chantier ?? -> fait\",\n  \"Numero de matricule supprime toujours les 0 -> fait\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Normalize whitespace (collapse runs of spaces/newlines) and curly\n// apostrophes so the comparison is robust to trailing spaces introduced by\n// multiple runs per paragraph and to quote-character variants.\nfunction normalize(s) {\n  return s.replace(/\\u2019/g, \"'\").replace(/\\s+/g, \" \").trim();\n}\n\nconst targets = textsToRemove.map(normalize);\n\nfor (let i = paragraphs.items.length - 1; i >= 0; i--) {\n  const para = paragraphs.items[i];\n  const normalized = normalize(para.text);\n  if (targets.includes(normalized)) {\n    para.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the resolved bug-tracker paragraphs from the document body.\n# These are the exact paragraph texts that must be deleted entirely\n# (whole paragraph, including every run inside it).\n\n$d = $word.ActiveDocument\n\n$rsquo = [char]0x2019\n\n$textsToRemove = @(\n  \"Matricule doit etre string ou int : parce que si int 00012 devient 12 -> fait\",\n  \"Mettre automatique la barre / pour les dates ?? -> fait\",\n  \"Modifier les vue pour ne pas utilise dispose sur la vue\",\n  (\"Est ce qu\" + $rsquo + \"on supprime tous les objet qui vont avec un genre on supprime les livraison d\" + $rsquo + \"un chantier si on supprime le chantier ?? -> fait\"),\n  \"Numero de matricule supprime toujours les 0 -> fait\"\n)\n\nfunction Normalize-Text($s) {\n  $s = $s -replace [char]0x2019, \"'\"\n  $s = $s -replace '\\s+', ' '\n  return $s.Trim()\n}\n\n$targets = @()\nforeach ($t in $textsToRemove) { $targets += (Normalize-Text $t) }\n\n# Walk paragraphs back-to-front so deleting one doesn't shift the indices\n# of the ones we still need to inspect.\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 1; $i--) {\n  $p = $d.Paragraphs.Item($i)\n  $ptext = Normalize-Text $p.Range.Text\n  if ($targets -contains $ptext) {\n    $p.Range.Delete()\n  }\n}\n"}
